# census-template.xlsx: add "Gender" and "Employee Number" columns (closes #419)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Gender" column after "Suffix" (new column E) ---
$ws.Columns("E:E").Insert()
$ws.Range("E1").Value = "Gender"
$ws.Columns("E:E").ColumnWidth = 6

# --- Insert "Employee Number" column after "Social Security Number" (new column O) ---
$ws.Columns("O:O").Insert()
$ws.Range("O1").Value = "Employee Number"
$ws.Columns("O:O").ColumnWidth = 15.5

# --- Give the "Date of Birth" data cell (K2) a date number format ---
$ws.Range("K2").NumberFormat = "mm-dd-yy"

# --- Update the active selection to match the saved view ---
[void]$ws.Range("P7").Select()
